# Reorders the names in the "Recorded By" column (column G) of the
# attendance/session-analysis sheet.
#
# Rule (derived from the target diff):
#   - Split the cell's comma-separated value into trimmed parts.
#   - If one of the parts is exactly "System" (case-sensitive), move it to
#     the front of the list, keeping the remaining parts in their original
#     order.
#   - Otherwise (no exact "System" entry), reverse the order of the parts.
#   - Cells with a single value are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7 ("Recorded By")
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }
    if (-not ($val -is [string])) {
        continue
    }
    if ($val.IndexOf(",") -lt 0) {
        continue
    }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if ($parts.Count -le 1) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
    } else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newVal = $newParts -join ", "

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
